# Insert a new weekly price record as row 106, pushing the existing
# rows 106-162 down to 107-163 (dimension grows from A1:R162 to A1:R163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 106..162 down by one row.
$ws.Rows.Item(106).Insert()

# Match the date-column number format used by the rest of column D.
$ws.Cells.Item(106, 4).NumberFormat = $ws.Cells.Item(107, 4).NumberFormat

# Populate the newly inserted row 106 with the new record.
$ws.Cells.Item(106, 1).Value2 = 4
$ws.Cells.Item(106, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(106, 3).Value2 = "Los Lagos"
$ws.Cells.Item(106, 4).Value2 = 44452
$ws.Cells.Item(106, 5).Value2 = 10
$ws.Cells.Item(106, 6).Value2 = 100112040
$ws.Cells.Item(106, 7).Value2 = "Cilantro"
$ws.Cells.Item(106, 8).Value2 = "Sin especificar"
$ws.Cells.Item(106, 9).Value2 = "Primera"
$ws.Cells.Item(106, 10).Value2 = 150
$ws.Cells.Item(106, 11).Value2 = 14000
$ws.Cells.Item(106, 12).Value2 = 14000
$ws.Cells.Item(106, 13).Value2 = 14000
$ws.Cells.Item(106, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(106, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(106, 16).Value2 = 389
$ws.Cells.Item(106, 17).Value2 = 36
$ws.Cells.Item(106, 18).Value2 = "Hortaliza"
